$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:G2").NumberFormat = "@"
$ws.Range("D2").Value = "332.42"
$ws.Range("E2").Value = "2.10%"
$ws.Range("F2").Value = "8-2-2023"
$ws.Range("G2").Value = "0"

$ws.Range("D3:G3").NumberFormat = "@"
$ws.Range("D3").Value = "45.60"
$ws.Range("E3").Value = "3.80%"
$ws.Range("F3").Value = "8-2-2023"
$ws.Range("G3").Value = "0"

$ws.Range("D4:G4").NumberFormat = "@"
$ws.Range("D4").Value = "5.585"
$ws.Range("E4").Value = "0.87%"
$ws.Range("F4").Value = "8-2-2023"
$ws.Range("G4").Value = "0"

$ws.Range("D5:G5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08320"
$ws.Range("E5").Value = "3.71%"
$ws.Range("F5").Value = "8-2-2023"
$ws.Range("G5").Value = "0"

$ws.Range("D6:G6").NumberFormat = "@"
$ws.Range("D6").Value = "2.039"
$ws.Range("E6").Value = "4.68%"
$ws.Range("F6").Value = "8-2-2023"
$ws.Range("G6").Value = "0"

$ws.Range("D7:G7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9817"
$ws.Range("E7").Value = "4.08%"
$ws.Range("F7").Value = "8-2-2023"
$ws.Range("G7").Value = "0"

$ws.Range("D8:G8").NumberFormat = "@"
$ws.Range("D8").Value = "0.1139"
$ws.Range("E8").Value = "2.34%"
$ws.Range("F8").Value = "8-2-2023"
$ws.Range("G8").Value = "0"

$ws.Range("D9:G9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1959"
$ws.Range("E9").Value = "6.86%"
$ws.Range("F9").Value = "8-2-2023"
$ws.Range("G9").Value = "0"

$ws.Range("D10:G10").NumberFormat = "@"
$ws.Range("D10").Value = "10.35"
$ws.Range("E10").Value = "-13.15%"
$ws.Range("F10").Value = "8-2-2023"
$ws.Range("G10").Value = "0"

$ws.Range("D11:G11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1008"
$ws.Range("E11").Value = "4.70%"
$ws.Range("F11").Value = "8-2-2023"
$ws.Range("G11").Value = "0"

$ws.Range("E12:G12").NumberFormat = "@"
$ws.Range("E12").Value = "-2.51%"
$ws.Range("F12").Value = "8-2-2023"
$ws.Range("G12").Value = "0"

$ws.Range("D13:G13").NumberFormat = "@"
$ws.Range("D13").Value = "0.1060"
$ws.Range("E13").Value = "-0.68%"
$ws.Range("F13").Value = "8-2-2023"
$ws.Range("G13").Value = "0"

$ws.Range("D14:G14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001259"
$ws.Range("E14").Value = "-0.61%"
$ws.Range("F14").Value = "8-2-2023"
$ws.Range("G14").Value = "0"

$ws.Range("D15:G15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005936"
$ws.Range("E15").Value = "3.79%"
$ws.Range("F15").Value = "8-2-2023"
$ws.Range("G15").Value = "0"

$ws.Range("D16:G16").NumberFormat = "@"
$ws.Range("D16").Value = "3.370"
$ws.Range("E16").Value = "-0.12%"
$ws.Range("F16").Value = "8-2-2023"
$ws.Range("G16").Value = "0"

$ws.Range("D17:G17").NumberFormat = "@"
$ws.Range("D17").Value = "4.439"
$ws.Range("E17").Value = "3.19%"
$ws.Range("F17").Value = "8-2-2023"
$ws.Range("G17").Value = "0"

$ws.Range("D18:G18").NumberFormat = "@"
$ws.Range("D18").Value = "2.616"
$ws.Range("E18").Value = "3.34%"
$ws.Range("F18").Value = "8-2-2023"
$ws.Range("G18").Value = "0"

$ws.Range("D19:G19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3333"
$ws.Range("E19").Value = "-4.14%"
$ws.Range("F19").Value = "8-2-2023"
$ws.Range("G19").Value = "0"

$ws.Range("D20:G20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1386"
$ws.Range("E20").Value = "-0.88%"
$ws.Range("F20").Value = "8-2-2023"
$ws.Range("G20").Value = "0"

$ws.Range("D21:G21").NumberFormat = "@"
$ws.Range("D21").Value = "0.2492"
$ws.Range("E21").Value = "-2.15%"
$ws.Range("F21").Value = "8-2-2023"
$ws.Range("G21").Value = "0"

$ws.Range("D22:G22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04112"
$ws.Range("E22").Value = "2.04%"
$ws.Range("F22").Value = "8-2-2023"
$ws.Range("G22").Value = "0"

$ws.Range("D23:G23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001301"
$ws.Range("E23").Value = "4.90%"
$ws.Range("F23").Value = "8-2-2023"
$ws.Range("G23").Value = "0"

$ws.Range("D24:G24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004428"
$ws.Range("E24").Value = "2.88%"
$ws.Range("F24").Value = "8-2-2023"
$ws.Range("G24").Value = "0"

$ws.Range("D25:G25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001279"
$ws.Range("E25").Value = "7.30%"
$ws.Range("F25").Value = "8-2-2023"
$ws.Range("G25").Value = "0"

$ws.Range("D26:G26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0003743"
$ws.Range("E26").Value = "-0.09%"
$ws.Range("F26").Value = "8-2-2023"
$ws.Range("G26").Value = "0"

$ws.Range("F27:G27").NumberFormat = "@"
$ws.Range("F27").Value = "8-2-2023"
$ws.Range("G27").Value = "0"

$ws.Range("F28:G28").NumberFormat = "@"
$ws.Range("F28").Value = "8-2-2023"
$ws.Range("G28").Value = "0"

$ws.Range("F29:G29").NumberFormat = "@"
$ws.Range("F29").Value = "8-2-2023"
$ws.Range("G29").Value = "0"

$ws.Range("F30:G30").NumberFormat = "@"
$ws.Range("F30").Value = "8-2-2023"
$ws.Range("G30").Value = "0"

$ws.Range("F31:G31").NumberFormat = "@"
$ws.Range("F31").Value = "8-2-2023"
$ws.Range("G31").Value = "0"

$ws.Range("F32:G32").NumberFormat = "@"
$ws.Range("F32").Value = "8-2-2023"
$ws.Range("G32").Value = "0"

$ws.Range("F33:G33").NumberFormat = "@"
$ws.Range("F33").Value = "8-2-2023"
$ws.Range("G33").Value = "0"

$ws.Range("F34:G34").NumberFormat = "@"
$ws.Range("F34").Value = "8-2-2023"
$ws.Range("G34").Value = "0"

$ws.Range("F35:G35").NumberFormat = "@"
$ws.Range("F35").Value = "8-2-2023"
$ws.Range("G35").Value = "0"

$ws.Range("F36:G36").NumberFormat = "@"
$ws.Range("F36").Value = "8-2-2023"
$ws.Range("G36").Value = "0"

$ws.Range("F37:G37").NumberFormat = "@"
$ws.Range("F37").Value = "8-2-2023"
$ws.Range("G37").Value = "0"

$ws.Range("D38:G38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02849"
$ws.Range("E38").Value = "13.03%"
$ws.Range("F38").Value = "8-2-2023"
$ws.Range("G38").Value = "0"

$ws.Range("D39:G39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05794"
$ws.Range("E39").Value = "5.31%"
$ws.Range("F39").Value = "8-2-2023"
$ws.Range("G39").Value = "0"

$ws.Range("D40:G40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007655"
$ws.Range("E40").Value = "1.73%"
$ws.Range("F40").Value = "8-2-2023"
$ws.Range("G40").Value = "0"

$ws.Range("D41:G41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1434"
$ws.Range("E41").Value = "3.54%"
$ws.Range("F41").Value = "8-2-2023"
$ws.Range("G41").Value = "0"

$ws.Range("D42:G42").NumberFormat = "@"
$ws.Range("D42").Value = "0.007704"
$ws.Range("E42").Value = "3.97%"
$ws.Range("F42").Value = "8-2-2023"
$ws.Range("G42").Value = "0"

$ws.Range("E43:G43").NumberFormat = "@"
$ws.Range("E43").Value = "-2.25%"
$ws.Range("F43").Value = "8-2-2023"
$ws.Range("G43").Value = "0"

$ws.Range("D44:G44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008129"
$ws.Range("E44").Value = "-2.26%"
$ws.Range("F44").Value = "8-2-2023"
$ws.Range("G44").Value = "0"

$ws.Range("D45:G45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00007205"
$ws.Range("E45").Value = "1.45%"
$ws.Range("F45").Value = "8-2-2023"
$ws.Range("G45").Value = "0"

$ws.Range("D46:G46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000751"
$ws.Range("E46").Value = "-0.10%"
$ws.Range("F46").Value = "8-2-2023"
$ws.Range("G46").Value = "0"

$ws.Range("D47:G47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0005805"
$ws.Range("E47").Value = "-0.11%"
$ws.Range("F47").Value = "8-2-2023"
$ws.Range("G47").Value = "0"

$ws.Range("D48:G48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003485"
$ws.Range("E48").Value = "-1.43%"
$ws.Range("F48").Value = "8-2-2023"
$ws.Range("G48").Value = "0"

$ws.Range("F49:G49").NumberFormat = "@"
$ws.Range("F49").Value = "8-2-2023"
$ws.Range("G49").Value = "0"

$ws.Range("D50:G50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002102"
$ws.Range("E50").Value = "-0.10%"
$ws.Range("F50").Value = "8-2-2023"
$ws.Range("G50").Value = "0"

$ws.Range("D51:G51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002002"
$ws.Range("E51").Value = "-0.10%"
$ws.Range("F51").Value = "8-2-2023"
$ws.Range("G51").Value = "0"
